# "duplicate container data merged"
#
# Sheet "Stock Report": row 14 was a duplicate of the container that also
# appears (merged) elsewhere; its row is cleared out except for the damage
# description, which moves into column X. The SL# numbering (column A) of
# every subsequent row shifts down by one. Two unrelated date cells (M6 and
# U6) pick up the date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# --- M6 / U6: apply the date number format (style index 4 already exists
#     in the workbook, so this reuses it rather than creating a new one) ---
$ws.Cells.Item(6, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(6, 21).NumberFormat = "YYYY-MM-DD"

# --- Row 14: clear every cell's content (keeping each cell's own style),
#     except column X which keeps the surviving damage description ---
$row = 14
for ($col = 1; $col -le 23; $col++) {
    $ws.Cells.Item($row, $col).ClearContents()
}
$ws.Cells.Item($row, 24).Value = "INTERIOR PANEL INK"
for ($col = 25; $col -le 28; $col++) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# --- Rows 15-24: SL# (column A) decrements by one now that the duplicate
#     row 14 entry no longer counts ---
for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}
